$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (s=1) from BC1 across the new header range BD1:BR1
$ws.Range("BC1").Copy($ws.Range("BD1:BR1"))

# Set new header labels for row 1
$ws.Range("BD1").Value = 'RV_ Tarjetas credito vigentes otros'
$ws.Range("BE1").Value = 'FV+F_ Tarjetas credito vigentes otros'
$ws.Range("BF1").Value = 'Threeshold_ Tarjetas credito vigentes otros'
$ws.Range("BG1").Value = 'RV_Numero de operaciones realizadas con tarjetas de crédito'
$ws.Range("BH1").Value = 'FV+F_Numero de operaciones realizadas con tarjetas de crédito'
$ws.Range("BI1").Value = 'Threeshold_Numero de operaciones realizadas con tarjetas de crédito'
$ws.Range("BJ1").Value = 'RV_Tarjetas vigentes'
$ws.Range("BK1").Value = 'FV+F_Tarjetas vigentes'
$ws.Range("BL1").Value = 'Threeshold_Tarjetas vigentes'
$ws.Range("BM1").Value = 'RV_Tarjetas vigentes VISA'
$ws.Range("BN1").Value = 'FV+F_Tarjetas vigentes VISA'
$ws.Range("BO1").Value = 'Threeshold_Tarjetas vigentes VISA'
$ws.Range("BP1").Value = 'RV_Tarjetas vigentes MASTERCARD'
$ws.Range("BQ1").Value = 'FV+F_Tarjetas vigentes MASTERCARD'
$ws.Range("BR1").Value = 'Threeshold_Tarjetas vigentes MASTERCARD'

# Fill data rows 2-27: RV/FV+F columns stay blank, Threeshold columns = 54
For ($r = 2; $r -le 27; $r++) {
    $ws.Range("BF$r").Value = 54
    $ws.Range("BI$r").Value = 54
    $ws.Range("BL$r").Value = 54
    $ws.Range("BO$r").Value = 54
    $ws.Range("BR$r").Value = 54
}
